# Shift every schedule date in column A (rows 2-11) forward by one day on
# every worksheet (Table 1 AM/PM .. Table 4 AM/PM). The last row (A11) ends
# up with its cell formatting cleared, matching the source edit.

$wb = $excel.ActiveWorkbook

$dates = @(
    "20/05/2024",
    "21/05/2024",
    "22/05/2024",
    "23/05/2024",
    "24/05/2024",
    "25/05/2024",
    "26/05/2024",
    "27/05/2024",
    "28/05/2024",
    "29/05/2024"
)

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $dates.Length; $i++) {
        $row = $i + 2
        $ws.Range("A$row").Value = $dates[$i]
    }
    # The final row's date cell loses its style (observed in the source edit).
    $ws.Range("A11").ClearFormats()
}
